$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Disable "Calculate before Save" (calcPr calcOnSave="0")
$excel.CalculateBeforeSave = $false

# Row 19 (Eagle Schematic / ERC): update notes text
$ws.Range("D19").Value = "Matched Guidelines. ERC Cleared. "

# Row 18 (Eagle Library): add new notes text
$ws.Range("D18").Value = "Matched Guidelines. "

# Row 20 (Eagle Layout / DRC): update notes text
$ws.Range("D20").Value = "Matched Guidelines. DRC is cleared"

# Update the active selection to D20
$ws.Range("D20").Select() | Out-Null
